$wb = $excel.ActiveWorkbook

# Add a new worksheet "ForeignTest" after the last existing sheet (GroupedItemTest)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "ForeignTest"

# Row 1 - label
$ws.Range("A1").Value = "C9"

# Row 9 - header row
$ws.Range("C9").Value = "Id"
$ws.Range("D9").Value = "TargetTestId"
$ws.Range("E9").Value = "Value"
$ws.Range("F9").Value = "StudentId"
$ws.Range("G9").Value = "비고"

# Row 10
$ws.Range("C10").Value = 1001
$ws.Range("D10").Value = 100
$ws.Range("E10").Value = "AAA"
$ws.Range("F10").Value = 20220001
$ws.Range("G10").Value = "학생이 있을수도 있고,"

# Row 11
$ws.Range("C11").Value = 1002
$ws.Range("D11").Value = 102
$ws.Range("E11").Value = "BBB"
$ws.Range("G11").Value = "학생이 없을수도 있습니다."

# Row 12
$ws.Range("C12").Value = 1003
$ws.Range("D12").Value = 104
$ws.Range("E12").Value = "CCC"
$ws.Range("F12").Value = 20220002

# Match the selection recorded for the new active sheet
$ws.Range("H22").Select() | Out-Null
